$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I (9) rows 2-10 held "deuteron" -- change the label to "d"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "d"
}

# Header row (A1:K1) becomes bold + centered
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Update the active selection to match the saved view
$null = $ws.Range("E21").Select()
